# Bump the "想去人数" (F column) counts by 1 for specific rows across sheets.
$wb = $excel.ActiveWorkbook

# Map of worksheet index (1-based) -> list of (cell, newValue)
$updates = @{
    1 = @(
        @{ Cell = "F10"; Value = 745 },
        @{ Cell = "F19"; Value = 68 },
        @{ Cell = "F23"; Value = 350 }
    )
    2 = @(
        @{ Cell = "F6";  Value = 471 },
        @{ Cell = "F19"; Value = 471 }
    )
    3 = @(
        @{ Cell = "F2"; Value = 1736 },
        @{ Cell = "F6"; Value = 2175 },
        @{ Cell = "F8"; Value = 833 }
    )
    4 = @(
        @{ Cell = "F2";  Value = 1736 },
        @{ Cell = "F4";  Value = 2175 },
        @{ Cell = "F9";  Value = 833 },
        @{ Cell = "F15"; Value = 745 },
        @{ Cell = "F24"; Value = 471 },
        @{ Cell = "F27"; Value = 68 },
        @{ Cell = "F31"; Value = 350 }
    )
}

foreach ($sheetIndex in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetIndex)
    foreach ($entry in $updates[$sheetIndex]) {
        $ws.Range($entry.Cell).Value = $entry.Value
    }
}
